$wb = $excel.ActiveWorkbook

# The two sheets that carry event data ("展览" and "全部类型") both need
# the same three "想去人数" (want-to-go count) updates applied to rows 2-4.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 66
    $ws.Range("F3").Value = 1304
    $ws.Range("F4").Value = 6
}
